$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix the border colors used by the tables "last row" style (currently row 96). ---
# The right/top thin-border edges used on the closing row move from pure white (FFFFFFFF)
# to an off-white (FFF8F9FA). Doing this on the *current* last row first means row 96s
# formatting can then be cloned onto the new, real last row (109) further down.
$lastRow = $ws.Range("A96:AA96")
$lastRow.Borders.Item(8).Color = 16447992
$lastRow.Borders.Item(11).Color = 16447992

# --- Step 2: capture that (now color-corrected) last-row format, for the new last row (109). ---
$ws.Range("A96:AA96").Copy()
$ws.Range("A109:AA109").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: re-format row 96 to the regular (non-last) alternating style, ---
# matching row 94 (same parity: the 13/14/15 style set).
$ws.Range("A94:AA94").Copy()
$ws.Range("A96:AA96").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: format the new rows 97-108, alternating like the existing data rows. ---
# Odd offset rows (97,99,101,103,105,107) take the "10/11/12" style (like row 95).
# Even offset rows (98,100,102,104,106,108) take the "13/14/15" style (like row 94).
$ws.Range("A95:AA95").Copy()
$ws.Range("A97:AA97").PasteSpecial(-4122)
$ws.Range("A99:AA99").PasteSpecial(-4122)
$ws.Range("A101:AA101").PasteSpecial(-4122)
$ws.Range("A103:AA103").PasteSpecial(-4122)
$ws.Range("A105:AA105").PasteSpecial(-4122)
$ws.Range("A107:AA107").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A94:AA94").Copy()
$ws.Range("A98:AA98").PasteSpecial(-4122)
$ws.Range("A100:AA100").PasteSpecial(-4122)
$ws.Range("A102:AA102").PasteSpecial(-4122)
$ws.Range("A104:AA104").PasteSpecial(-4122)
$ws.Range("A106:AA106").PasteSpecial(-4122)
$ws.Range("A108:AA108").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 5: fill in cell values for the 13 new survey responses (rows 97-109). ---
# Row 97
$ws.Cells.Item(97, 1).Value = 45608.67305858796
$ws.Cells.Item(97, 2).Value = "minseok1937@gmail.com"
$ws.Cells.Item(97, 3).Value = "경영"
$ws.Cells.Item(97, 4).Value = 20202915
$ws.Cells.Item(97, 5).Value = "김민석"
$ws.Cells.Item(97, 6).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 7).Value = "2. 별로 아니다"
$ws.Cells.Item(97, 8).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 9).Value = "5. 매우 그렇다"
$ws.Cells.Item(97, 10).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(97, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 14).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 15).Value = "2. 별로 아니다"
$ws.Cells.Item(97, 16).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 17).Value = "1. 전혀 아니다"
$ws.Cells.Item(97, 18).Value = 3.0
$ws.Cells.Item(97, 19).Value = 3.0
$ws.Cells.Item(97, 20).Value = 10.0
$ws.Cells.Item(97, 21).Value = 4.0
$ws.Cells.Item(97, 22).Value = 3.0
$ws.Cells.Item(97, 23).Value = "낮음"
$ws.Cells.Item(97, 24).Value = "낮음"
$ws.Cells.Item(97, 25).Value = "높음"
$ws.Cells.Item(97, 26).Value = "낮음"
$ws.Cells.Item(97, 27).Value = "낮음"
# Row 98
$ws.Cells.Item(98, 1).Value = 45608.675921076385
$ws.Cells.Item(98, 2).Value = "h20191240@glab.hallym.ac.kr"
$ws.Cells.Item(98, 3).Value = "영어영문학과"
$ws.Cells.Item(98, 4).Value = 20191240
$ws.Cells.Item(98, 5).Value = "홍이래"
$ws.Cells.Item(98, 6).Value = "1. 전혀 아니다"
$ws.Cells.Item(98, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(98, 8).Value = "1. 전혀 아니다"
$ws.Cells.Item(98, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(98, 10).Value = "2. 별로 아니다"
$ws.Cells.Item(98, 11).Value = "4. 약간 그렇다"
$ws.Cells.Item(98, 12).Value = "2. 별로 아니다"
$ws.Cells.Item(98, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(98, 14).Value = "1. 전혀 아니다"
$ws.Cells.Item(98, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(98, 16).Value = "1. 전혀 아니다"
$ws.Cells.Item(98, 17).Value = "5. 매우 그렇다"
$ws.Cells.Item(98, 18).Value = 5.0
$ws.Cells.Item(98, 19).Value = 6.0
$ws.Cells.Item(98, 20).Value = 9.0
$ws.Cells.Item(98, 21).Value = 14.0
$ws.Cells.Item(98, 22).Value = 7.0
$ws.Cells.Item(98, 23).Value = "중하"
$ws.Cells.Item(98, 24).Value = "중하"
$ws.Cells.Item(98, 25).Value = "높음"
$ws.Cells.Item(98, 26).Value = "높음"
$ws.Cells.Item(98, 27).Value = "중하"
# Row 99
$ws.Cells.Item(99, 1).Value = 45608.685173252314
$ws.Cells.Item(99, 2).Value = "seollo020531@naver.com"
$ws.Cells.Item(99, 3).Value = "소프트웨어학과"
$ws.Cells.Item(99, 4).Value = 20225175
$ws.Cells.Item(99, 5).Value = "설창원"
$ws.Cells.Item(99, 6).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 7).Value = "3. 중간이다"
$ws.Cells.Item(99, 8).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(99, 10).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(99, 13).Value = "4. 약간 그렇다"
$ws.Cells.Item(99, 14).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 15).Value = "3. 중간이다"
$ws.Cells.Item(99, 16).Value = "2. 별로 아니다"
$ws.Cells.Item(99, 17).Value = "4. 약간 그렇다"
$ws.Cells.Item(99, 18).Value = 4.0
$ws.Cells.Item(99, 19).Value = 5.0
$ws.Cells.Item(99, 20).Value = 8.0
$ws.Cells.Item(99, 21).Value = 12.0
$ws.Cells.Item(99, 22).Value = 8.0
$ws.Cells.Item(99, 23).Value = "낮음"
$ws.Cells.Item(99, 24).Value = "중하"
$ws.Cells.Item(99, 25).Value = "중상"
$ws.Cells.Item(99, 26).Value = "높음"
$ws.Cells.Item(99, 27).Value = "중상"
# Row 100
$ws.Cells.Item(100, 1).Value = 45608.69253414352
$ws.Cells.Item(100, 2).Value = "algus5661@naver.com"
$ws.Cells.Item(100, 3).Value = "사회복지학부"
$ws.Cells.Item(100, 4).Value = 20242355
$ws.Cells.Item(100, 5).Value = "최미현"
$ws.Cells.Item(100, 6).Value = "3. 중간이다"
$ws.Cells.Item(100, 7).Value = "3. 중간이다"
$ws.Cells.Item(100, 8).Value = "1. 전혀 아니다"
$ws.Cells.Item(100, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(100, 10).Value = "2. 별로 아니다"
$ws.Cells.Item(100, 11).Value = "4. 약간 그렇다"
$ws.Cells.Item(100, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(100, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(100, 14).Value = "1. 전혀 아니다"
$ws.Cells.Item(100, 15).Value = "3. 중간이다"
$ws.Cells.Item(100, 16).Value = "3. 중간이다"
$ws.Cells.Item(100, 17).Value = "5. 매우 그렇다"
$ws.Cells.Item(100, 18).Value = 7.0
$ws.Cells.Item(100, 19).Value = 6.0
$ws.Cells.Item(100, 20).Value = 6.0
$ws.Cells.Item(100, 21).Value = 9.0
$ws.Cells.Item(100, 22).Value = 5.0
$ws.Cells.Item(100, 23).Value = "중상"
$ws.Cells.Item(100, 24).Value = "중하"
$ws.Cells.Item(100, 25).Value = "중하"
$ws.Cells.Item(100, 26).Value = "낮음"
$ws.Cells.Item(100, 27).Value = "낮음"
# Row 101
$ws.Cells.Item(101, 1).Value = 45608.69424950231
$ws.Cells.Item(101, 2).Value = "oepdwrtyy@gmail.com"
$ws.Cells.Item(101, 3).Value = "인문학부"
$ws.Cells.Item(101, 4).Value = 20241003
$ws.Cells.Item(101, 5).Value = "강종현"
$ws.Cells.Item(101, 6).Value = "2. 별로 아니다"
$ws.Cells.Item(101, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(101, 8).Value = "5. 매우 그렇다"
$ws.Cells.Item(101, 9).Value = "3. 중간이다"
$ws.Cells.Item(101, 10).Value = "4. 약간 그렇다"
$ws.Cells.Item(101, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(101, 12).Value = "2. 별로 아니다"
$ws.Cells.Item(101, 13).Value = "4. 약간 그렇다"
$ws.Cells.Item(101, 14).Value = "2. 별로 아니다"
$ws.Cells.Item(101, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(101, 16).Value = "2. 별로 아니다"
$ws.Cells.Item(101, 17).Value = "3. 중간이다"
$ws.Cells.Item(101, 18).Value = 4.0
$ws.Cells.Item(101, 19).Value = 8.0
$ws.Cells.Item(101, 20).Value = 7.0
$ws.Cells.Item(101, 21).Value = 11.0
$ws.Cells.Item(101, 22).Value = 11.0
$ws.Cells.Item(101, 23).Value = "낮음"
$ws.Cells.Item(101, 24).Value = "중상"
$ws.Cells.Item(101, 25).Value = "중상"
$ws.Cells.Item(101, 26).Value = "중하"
$ws.Cells.Item(101, 27).Value = "중상"
# Row 102
$ws.Cells.Item(102, 1).Value = 45608.765146342594
$ws.Cells.Item(102, 2).Value = "eung4077@gmail.com"
$ws.Cells.Item(102, 3).Value = "사회학과"
$ws.Cells.Item(102, 4).Value = 20242201
$ws.Cells.Item(102, 5).Value = "강은결"
$ws.Cells.Item(102, 6).Value = "2. 별로 아니다"
$ws.Cells.Item(102, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 8).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 10).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(102, 12).Value = "2. 별로 아니다"
$ws.Cells.Item(102, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(102, 14).Value = "3. 중간이다"
$ws.Cells.Item(102, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 16).Value = "2. 별로 아니다"
$ws.Cells.Item(102, 17).Value = "4. 약간 그렇다"
$ws.Cells.Item(102, 18).Value = 4.0
$ws.Cells.Item(102, 19).Value = 8.0
$ws.Cells.Item(102, 20).Value = 7.0
$ws.Cells.Item(102, 21).Value = 12.0
$ws.Cells.Item(102, 22).Value = 7.0
$ws.Cells.Item(102, 23).Value = "낮음"
$ws.Cells.Item(102, 24).Value = "중상"
$ws.Cells.Item(102, 25).Value = "중상"
$ws.Cells.Item(102, 26).Value = "중하"
$ws.Cells.Item(102, 27).Value = "낮음"
# Row 103
$ws.Cells.Item(103, 1).Value = 45608.83104348379
$ws.Cells.Item(103, 2).Value = "syw050819@naver.com"
$ws.Cells.Item(103, 3).Value = "간호학과"
$ws.Cells.Item(103, 4).Value = 20246251
$ws.Cells.Item(103, 5).Value = "신예원"
$ws.Cells.Item(103, 6).Value = "4. 약간 그렇다"
$ws.Cells.Item(103, 7).Value = "2. 별로 아니다"
$ws.Cells.Item(103, 8).Value = "3. 중간이다"
$ws.Cells.Item(103, 9).Value = "3. 중간이다"
$ws.Cells.Item(103, 10).Value = "1. 전혀 아니다"
$ws.Cells.Item(103, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(103, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(103, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(103, 14).Value = "4. 약간 그렇다"
$ws.Cells.Item(103, 15).Value = "3. 중간이다"
$ws.Cells.Item(103, 16).Value = "2. 별로 아니다"
$ws.Cells.Item(103, 17).Value = "3. 중간이다"
$ws.Cells.Item(103, 18).Value = 6.0
$ws.Cells.Item(103, 19).Value = 4.0
$ws.Cells.Item(103, 20).Value = 5.0
$ws.Cells.Item(103, 21).Value = 10.0
$ws.Cells.Item(103, 22).Value = 6.0
$ws.Cells.Item(103, 23).Value = "중하"
$ws.Cells.Item(103, 24).Value = "낮음"
$ws.Cells.Item(103, 25).Value = "중하"
$ws.Cells.Item(103, 26).Value = "낮음"
$ws.Cells.Item(103, 27).Value = "낮음"
# Row 104
$ws.Cells.Item(104, 1).Value = 45608.83613487269
$ws.Cells.Item(104, 2).Value = "kimguswls6685@naver.com"
$ws.Cells.Item(104, 3).Value = "콘텐츠IT전공"
$ws.Cells.Item(104, 4).Value = 20215144
$ws.Cells.Item(104, 5).Value = "김현진"
$ws.Cells.Item(104, 6).Value = "2. 별로 아니다"
$ws.Cells.Item(104, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(104, 8).Value = "3. 중간이다"
$ws.Cells.Item(104, 9).Value = "3. 중간이다"
$ws.Cells.Item(104, 10).Value = "1. 전혀 아니다"
$ws.Cells.Item(104, 11).Value = "3. 중간이다"
$ws.Cells.Item(104, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(104, 13).Value = "2. 별로 아니다"
$ws.Cells.Item(104, 14).Value = "3. 중간이다"
$ws.Cells.Item(104, 15).Value = "2. 별로 아니다"
$ws.Cells.Item(104, 16).Value = "1. 전혀 아니다"
$ws.Cells.Item(104, 17).Value = "3. 중간이다"
$ws.Cells.Item(104, 18).Value = 5.0
$ws.Cells.Item(104, 19).Value = 3.0
$ws.Cells.Item(104, 20).Value = 6.0
$ws.Cells.Item(104, 21).Value = 12.0
$ws.Cells.Item(104, 22).Value = 8.0
$ws.Cells.Item(104, 23).Value = "중하"
$ws.Cells.Item(104, 24).Value = "낮음"
$ws.Cells.Item(104, 25).Value = "중하"
$ws.Cells.Item(104, 26).Value = "중상"
$ws.Cells.Item(104, 27).Value = "낮음"
# Row 105
$ws.Cells.Item(105, 1).Value = 45608.846135671294
$ws.Cells.Item(105, 2).Value = "ertyhx3@gmail.com"
$ws.Cells.Item(105, 3).Value = "광고홍보학과 "
$ws.Cells.Item(105, 4).Value = 20242607
$ws.Cells.Item(105, 5).Value = "김미소"
$ws.Cells.Item(105, 6).Value = "1. 전혀 아니다"
$ws.Cells.Item(105, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(105, 8).Value = "5. 매우 그렇다"
$ws.Cells.Item(105, 9).Value = "2. 별로 아니다"
$ws.Cells.Item(105, 10).Value = "4. 약간 그렇다"
$ws.Cells.Item(105, 11).Value = "2. 별로 아니다"
$ws.Cells.Item(105, 12).Value = "3. 중간이다"
$ws.Cells.Item(105, 13).Value = "4. 약간 그렇다"
$ws.Cells.Item(105, 14).Value = "3. 중간이다"
$ws.Cells.Item(105, 15).Value = "5. 매우 그렇다"
$ws.Cells.Item(105, 16).Value = "3. 중간이다"
$ws.Cells.Item(105, 17).Value = "5. 매우 그렇다"
$ws.Cells.Item(105, 18).Value = 3.0
$ws.Cells.Item(105, 19).Value = 4.0
$ws.Cells.Item(105, 20).Value = 3.0
$ws.Cells.Item(105, 21).Value = 3.0
$ws.Cells.Item(105, 22).Value = 4.0
$ws.Cells.Item(105, 23).Value = "중하"
$ws.Cells.Item(105, 24).Value = "중상"
$ws.Cells.Item(105, 25).Value = "중하"
$ws.Cells.Item(105, 26).Value = "중하"
$ws.Cells.Item(105, 27).Value = "중상"
# Row 106
$ws.Cells.Item(106, 1).Value = 45608.849338796295
$ws.Cells.Item(106, 2).Value = "withhowon@gmail.com"
$ws.Cells.Item(106, 3).Value = "간호학과"
$ws.Cells.Item(106, 4).Value = 20246245
$ws.Cells.Item(106, 5).Value = "서호원"
$ws.Cells.Item(106, 6).Value = "4. 약간 그렇다"
$ws.Cells.Item(106, 7).Value = "3. 중간이다"
$ws.Cells.Item(106, 8).Value = "2. 별로 아니다"
$ws.Cells.Item(106, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(106, 10).Value = "4. 약간 그렇다"
$ws.Cells.Item(106, 11).Value = "3. 중간이다"
$ws.Cells.Item(106, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(106, 13).Value = "2. 별로 아니다"
$ws.Cells.Item(106, 14).Value = "2. 별로 아니다"
$ws.Cells.Item(106, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(106, 16).Value = "3. 중간이다"
$ws.Cells.Item(106, 17).Value = "5. 매우 그렇다"
$ws.Cells.Item(106, 18).Value = 7.0
$ws.Cells.Item(106, 19).Value = 8.0
$ws.Cells.Item(106, 20).Value = 8.0
$ws.Cells.Item(106, 21).Value = 13.0
$ws.Cells.Item(106, 22).Value = 9.0
$ws.Cells.Item(106, 23).Value = "중상"
$ws.Cells.Item(106, 24).Value = "중상"
$ws.Cells.Item(106, 25).Value = "중상"
$ws.Cells.Item(106, 26).Value = "중상"
$ws.Cells.Item(106, 27).Value = "중하"
# Row 107
$ws.Cells.Item(107, 1).Value = 45608.85875461806
$ws.Cells.Item(107, 2).Value = "snp040609@naver.com"
$ws.Cells.Item(107, 3).Value = "경영학과"
$ws.Cells.Item(107, 4).Value = 20242957
$ws.Cells.Item(107, 5).Value = "박세나"
$ws.Cells.Item(107, 6).Value = "4. 약간 그렇다"
$ws.Cells.Item(107, 7).Value = "4. 약간 그렇다"
$ws.Cells.Item(107, 8).Value = "2. 별로 아니다"
$ws.Cells.Item(107, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(107, 10).Value = "3. 중간이다"
$ws.Cells.Item(107, 11).Value = "5. 매우 그렇다"
$ws.Cells.Item(107, 12).Value = "2. 별로 아니다"
$ws.Cells.Item(107, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(107, 14).Value = "2. 별로 아니다"
$ws.Cells.Item(107, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(107, 16).Value = "3. 중간이다"
$ws.Cells.Item(107, 17).Value = "4. 약간 그렇다"
$ws.Cells.Item(107, 18).Value = 9.0
$ws.Cells.Item(107, 19).Value = 7.0
$ws.Cells.Item(107, 20).Value = 6.0
$ws.Cells.Item(107, 21).Value = 10.0
$ws.Cells.Item(107, 22).Value = 6.0
$ws.Cells.Item(107, 23).Value = "높음"
$ws.Cells.Item(107, 24).Value = "중상"
$ws.Cells.Item(107, 25).Value = "중하"
$ws.Cells.Item(107, 26).Value = "낮음"
$ws.Cells.Item(107, 27).Value = "낮음"
# Row 108
$ws.Cells.Item(108, 1).Value = 45608.86104516203
$ws.Cells.Item(108, 2).Value = "mt1661@naver.com"
$ws.Cells.Item(108, 3).Value = "콘탠츠IT전공"
$ws.Cells.Item(108, 4).Value = 20215239
$ws.Cells.Item(108, 5).Value = "정성민"
$ws.Cells.Item(108, 6).Value = "3. 중간이다"
$ws.Cells.Item(108, 7).Value = "2. 별로 아니다"
$ws.Cells.Item(108, 8).Value = "2. 별로 아니다"
$ws.Cells.Item(108, 9).Value = "2. 별로 아니다"
$ws.Cells.Item(108, 10).Value = "3. 중간이다"
$ws.Cells.Item(108, 11).Value = "3. 중간이다"
$ws.Cells.Item(108, 12).Value = "4. 약간 그렇다"
$ws.Cells.Item(108, 13).Value = "5. 매우 그렇다"
$ws.Cells.Item(108, 14).Value = "2. 별로 아니다"
$ws.Cells.Item(108, 15).Value = "5. 매우 그렇다"
$ws.Cells.Item(108, 16).Value = "3. 중간이다"
$ws.Cells.Item(108, 17).Value = "1. 전혀 아니다"
$ws.Cells.Item(108, 18).Value = 6.0
$ws.Cells.Item(108, 19).Value = 9.0
$ws.Cells.Item(108, 20).Value = 7.0
$ws.Cells.Item(108, 21).Value = 5.0
$ws.Cells.Item(108, 22).Value = 7.0
$ws.Cells.Item(108, 23).Value = "중하"
$ws.Cells.Item(108, 24).Value = "중상"
$ws.Cells.Item(108, 25).Value = "중상"
$ws.Cells.Item(108, 26).Value = "중하"
$ws.Cells.Item(108, 27).Value = "중상"
# Row 109
$ws.Cells.Item(109, 1).Value = 45608.862400868056
$ws.Cells.Item(109, 2).Value = "rer220@naver.com"
$ws.Cells.Item(109, 3).Value = "콘텐츠IT"
$ws.Cells.Item(109, 4).Value = 20205124
$ws.Cells.Item(109, 5).Value = "김대명"
$ws.Cells.Item(109, 6).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 7).Value = "2. 별로 아니다"
$ws.Cells.Item(109, 8).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 9).Value = "4. 약간 그렇다"
$ws.Cells.Item(109, 10).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 11).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 12).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 13).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 14).Value = "1. 전혀 아니다"
$ws.Cells.Item(109, 15).Value = "4. 약간 그렇다"
$ws.Cells.Item(109, 16).Value = "3. 중간이다"
$ws.Cells.Item(109, 17).Value = "4. 약간 그렇다"
$ws.Cells.Item(109, 18).Value = 2.0
$ws.Cells.Item(109, 19).Value = 5.0
$ws.Cells.Item(109, 20).Value = 5.0
$ws.Cells.Item(109, 21).Value = 7.0
$ws.Cells.Item(109, 22).Value = 5.0
$ws.Cells.Item(109, 23).Value = "낮음"
$ws.Cells.Item(109, 24).Value = "중하"
$ws.Cells.Item(109, 25).Value = "중하"
$ws.Cells.Item(109, 26).Value = "중상"
$ws.Cells.Item(109, 27).Value = "중하"

# --- Step 6: extend the structured table to cover the new rows. ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:AA109"))
